$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 20 (shifts all existing data, incl. rows
# 20-177, down by 2 -> new rows 22-179). This matches the diff: every
# existing row's content reappears two rows further down, the dimension
# grows from A1:R177 to A1:R179, and a brand-new "week" of data
# (Primera/Segunda) is added at the (now blank) rows 20 and 21.
$ws.Rows("20:21").Insert()

# --- Row 20: Primera ---
$ws.Cells.Item(20, 1).Value2 = 8
$ws.Cells.Item(20, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(20, 3).Value2 = "Coquimbo"
$ws.Cells.Item(20, 4).Value2 = 44490
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 5).Value2 = 4
$ws.Cells.Item(20, 6).Value2 = 100114014
$ws.Cells.Item(20, 7).Value2 = "Betarraga"
$ws.Cells.Item(20, 8).Value2 = "Sin especificar"
$ws.Cells.Item(20, 9).Value2 = "Primera"
$ws.Cells.Item(20, 10).Value2 = 3000
$ws.Cells.Item(20, 11).Value2 = 450
$ws.Cells.Item(20, 12).Value2 = 500
$ws.Cells.Item(20, 13).Value2 = 475
$ws.Cells.Item(20, 14).Value2 = "`$/paquete 3 unidades"
$ws.Cells.Item(20, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(20, 16).Value2 = 158
$ws.Cells.Item(20, 17).Value2 = 3
$ws.Cells.Item(20, 18).Value2 = "Hortaliza"

# --- Row 21: Segunda ---
$ws.Cells.Item(21, 1).Value2 = 8
$ws.Cells.Item(21, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(21, 3).Value2 = "Coquimbo"
$ws.Cells.Item(21, 4).Value2 = 44490
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 5).Value2 = 4
$ws.Cells.Item(21, 6).Value2 = 100114014
$ws.Cells.Item(21, 7).Value2 = "Betarraga"
$ws.Cells.Item(21, 8).Value2 = "Sin especificar"
$ws.Cells.Item(21, 9).Value2 = "Segunda"
$ws.Cells.Item(21, 10).Value2 = 1400
$ws.Cells.Item(21, 11).Value2 = 350
$ws.Cells.Item(21, 12).Value2 = 400
$ws.Cells.Item(21, 13).Value2 = 375
$ws.Cells.Item(21, 14).Value2 = "`$/paquete 3 unidades"
$ws.Cells.Item(21, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(21, 16).Value2 = 125
$ws.Cells.Item(21, 17).Value2 = 3
$ws.Cells.Item(21, 18).Value2 = "Hortaliza"
